# "testing full row 1"
# - Unmerge the A1:G1 title band
# - Fill B1:G1 with a "placeholder" label (previously empty, merge-only cells)
# - Explicitly size rows 1 & 2 to 12.75pt (matches the post-edit row heights)
# - Move the active selection to G1 (single cell, was C2:G2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The title band A1:G1 was merged; unmerge it so each cell can hold its own value.
[void]$ws.Range("A1:G1").UnMerge()

# Clear the inherited "merged blank" formatting from B1:G1 before writing values,
# then give each cell the placeholder text.
[void]$ws.Range("B1:G1").ClearFormats()
$ws.Range("B1").Value = "placeholder"
$ws.Range("C1").Value = "placeholder"
$ws.Range("D1").Value = "placeholder"
$ws.Range("E1").Value = "placeholder"
$ws.Range("F1").Value = "placeholder"
$ws.Range("G1").Value = "placeholder"

# Row heights settle to 12.75 once the band is no longer a single merged cell.
$ws.Rows.Item(1).RowHeight = 12.75
$ws.Rows.Item(2).RowHeight = 12.75

# Active cell moves to G1.
[void]$ws.Range("G1").Select()
